$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new dialog rows (sailor dock: fill food / no sailors available).
$ws.Range("A12").Value = "dialog_fill_food"
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = 7
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = "name_dock_sailor"
$ws.Range("F12").Value = 0

$ws.Range("A13").Value = "dialog_no_sailors"
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = "name_dock_sailor"
$ws.Range("F13").Value = 0

# Move the active selection to match the author's saved cursor position.
$ws.Range("F9").Select() | Out-Null
